$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet tracks analyst ratings week over week. Each new edit prepends
# two new "as of" date columns (this time Jun_15 and Jun_17) right after the
# rank column (A), pushing the existing date columns to the right. The new
# columns start out as copies of the existing placeholder column (B, "UN")
# since no rating data has been recorded for those dates yet.

# Insert the "Jun_15" column in position B (copy column B's placeholder
# values so every row gets "UN" to start with), shifting the existing
# columns (Jun_13, Jun_10) one column to the right.
$ws.Columns("B").Copy()
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "Jun_15"

# Insert the "Jun_17" column in position B as well (copy column B's
# placeholder values again). This pushes the "Jun_15" column (and
# everything right of it) one more column to the right, so the newest
# week ends up left-most.
$ws.Columns("B").Copy()
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "Jun_17"

# Match the original column C width (8.0) on the three now-adjacent columns
# C, D and E (the two new columns plus the shifted-over original column).
$ws.Columns("C").ColumnWidth = 7.1
$ws.Columns("D").ColumnWidth = 7.1
$ws.Columns("E").ColumnWidth = 7.1

Write-Host "Inserted Jun_15 and Jun_17 columns"
